# Generate Report for Handback
# Row 7 (the "cbe87332-5857-413b-9b79-7e0123e0647e" entry) on both the
# "zh-cn" and "de-de" status sheets now has a completed handback: a
# "Latest Target File" hyperlink, a "Latest Handback File"/"Latest Handback
# DateTime" pair, and an "Error Detail" message noting the handback file is
# stale.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3beabf9e9fec52863b80cddb7476b39e4266d98/e2e/cbe87332-5857-413b-9b79-7e0123e0647e.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc1d4ac8639ab57148b3a2857afae88c81f731fe/e2e/cbe87332-5857-413b-9b79-7e0123e0647e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3beabf9e9fec52863b80cddb7476b39e4266d98/e2e/cbe87332-5857-413b-9b79-7e0123e0647e.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", "cbe87332-5857-413b-9b79-7e0123e0647e.md")

$wsZh.Range("J7").Value = $wsZh.Range("G7").Text
$wsZh.Range("K7").Value = "2016-09-06 21:15:09"
$wsZh.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", "cbe87332-5857-413b-9b79-7e0123e0647e.md")

$wsDe.Range("J7").Value = $wsDe.Range("G7").Text
$wsDe.Range("K7").Value = "2016-09-06 21:15:27"
$wsDe.Range("P7").Value = $errorDetail
